$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 58, shifting existing rows 58..119 down to 59..120.
# Excel's Insert operation copies formatting from the row above (row 57),
# which matches the row style (date style on column D) seen in the target file.
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new weekly data point.
# Most fields mirror the (old) row 58 data, except Fecha (D) and Volumen (J).
$ws.Cells.Item(58, 1).Value = 7
$ws.Cells.Item(58, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(58, 3).Value = "Ñuble"
$ws.Cells.Item(58, 4).Value = [DateTime]"2023-09-06"
$ws.Cells.Item(58, 5).Value = 16
$ws.Cells.Item(58, 6).Value = 100112001
$ws.Cells.Item(58, 7).Value = "Berenjena"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 80
$ws.Cells.Item(58, 11).Value = 10000
$ws.Cells.Item(58, 12).Value = 10000
$ws.Cells.Item(58, 13).Value = 10000
$ws.Cells.Item(58, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(58, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(58, 16).Value = 167
$ws.Cells.Item(58, 17).Value = 60
$ws.Cells.Item(58, 18).Value = "Hortaliza"

# Ensure the date cell keeps the expected date/time number format used by the rest of column D.
$ws.Cells.Item(58, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
